# Add two new columns (I: "I0", J: "IF") to the sheet, mirroring the
# existing header style used by H1 ("IP") and filling in the numeric
# data for rows 2-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - set their text first, then copy the style from H1 so I1/J1
# match the other bold/centered header cells.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Data rows 2-16 for columns I (I0) and J (IF)
$data = @(
    @(8, 8),   # row 2
    @(8, 8),   # row 3
    @(8, 8),   # row 4
    @(8, 8),   # row 5
    @(9, 9),   # row 6
    @(2, 3),   # row 7
    @(6, 7),   # row 8
    @(6, 7),   # row 9
    @(8, 8),   # row 10
    @(5, 7),   # row 11
    @(8, 9),   # row 12
    @(5, 6),   # row 13
    @(7, 7),   # row 14
    @(3, 4),   # row 15
    @(8, 8)    # row 16
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
